$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 15: a hyperlinked display_id, the PCR attachment filename, and a
# CONCAT formula describing it.
$ws.Hyperlinks.Add($ws.Range("A15"), "https://synbiohub.org/user/jhay/Johnny190421/sl0199_flatten/1.0.0", "", "", "https://synbiohub.org/user/jhay/Johnny190421/sl0199_flatten/1.0.0")
$ws.Range("A15").Value = "sl0199_flatten"
$ws.Range("B15").Value = "flanks_short.xlsx"
$ws.Range("D15").Formula = "=CONCAT(B15, "" is a PCR file"")"

# Match the author's final selection on the new row.
[void]$ws.Range("A15").Select()
